# Auto-generated edit script applying cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.220.28"
$ws.Range("E2").Value = "  +7.88%  "

$ws.Range("D3").Value = "2.584.84"
$ws.Range("E3").Value = "  +10.10%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "503.75"
$ws.Range("E5").Value = "  +6.72%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.84"
$ws.Range("E6").Value = "  +8.77%  "

$ws.Range("E7").Value = "  +24.99%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.993"
$ws.Range("E8").Value = "  -0.62%  "

$ws.Range("D9").Value = "2.582.07"
$ws.Range("E9").Value = "  +10.05%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.13"
$ws.Range("E10").Value = "  +13.20%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.102"
$ws.Range("E11").Value = "  +6.89%  "

$ws.Range("E12").Value = "  +6.86%  "

$ws.Range("E13").Value = "  +1.79%  "

$ws.Range("D14").Value = "2.976.09"
$ws.Range("E14").Value = "  +7.92%  "

$ws.Range("D15").Value = "59.151.34"
$ws.Range("E15").Value = "  +7.77%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.69"
$ws.Range("E16").Value = "  +8.64%  "

$ws.Range("E17").Value = "  +5.54%  "

$ws.Range("D18").Value = "2.569.87"
$ws.Range("E18").Value = "  +9.20%  "

$ws.Range("E19").Value = "  +5.08%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "332.98"
$ws.Range("E20").Value = "  +7.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.33"
$ws.Range("E21").Value = "  +8.03%  "

$ws.Range("E22").Value = "  +7.96%  "

$ws.Range("E23").Value = "  +0.58%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "59.66"
$ws.Range("E24").Value = "  +6.82%  "

$ws.Range("E25").Value = "  +6.37%  "

$ws.Range("E26").Value = "  +8.68%  "

$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.994"
$ws.Range("E27").Value = "  -0.58%  "

$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.654.82"
$ws.Range("E28").Value = "  +8.43%  "

$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.40"
$ws.Range("E29").Value = "  +3.95%  "

$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0824"
$ws.Range("E30").Value = "  +10.18%  "

$ws.Range("E31").Value = "  -0.35%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "157.58"
$ws.Range("E32").Value = "  +7.10%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.29"
$ws.Range("E33").Value = "  +7.35%  "

$ws.Range("E34").Value = "  +7.12%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.50"
$ws.Range("E35").Value = "  +9.64%  "

$ws.Range("E36").Value = "  +9.84%  "

$ws.Range("E37").Value = "  +9.42%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.842"
$ws.Range("E38").Value = "  +3.17%  "

$ws.Range("E39").Value = "  +12.15%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.45"
$ws.Range("E40").Value = "  +8.39%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "35.18"
$ws.Range("E41").Value = "  +5.65%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "290.26"
$ws.Range("E42").Value = "  +15.54%  "

$ws.Range("E43").Value = "  +7.55%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.626"
$ws.Range("E44").Value = "  +8.94%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0563"
$ws.Range("E45").Value = "  +7.88%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.995"
$ws.Range("E46").Value = "  -0.28%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.765"
$ws.Range("E47").Value = "  +22.37%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.14"
$ws.Range("E48").Value = "  +14.65%  "

$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.84"
$ws.Range("E49").Value = "  +10.13%  "

$ws.Range("E50").Value = "  +7.05%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.25"
$ws.Range("E51").Value = "  +1.03%  "
